$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B width (closest achievable value on this engine's pixel grid
# to the target width of 15.42578125)
$ws.Columns.Item(2).ColumnWidth = 14.65

# Update cell values
$ws.Range("A1").Value = -0.022571533854122475
$ws.Range("B1").Value = -0.022620277743222433

$ws.Range("A2").Value = -0.022936736118728825
$ws.Range("B2").Value = -0.043120838565710057

$ws.Range("A3").Value = -0.0011949061696823679
$ws.Range("B3").Value = -0.0010070942439798769

$ws.Range("A4").Value = -0.079941757706758973
$ws.Range("B4").Value = -0.079976916014463731
